# Applies the "Anonimyzed fedcore" update:
#  - C1/D1 (sheet1) and C1/D1/F1/G1 (sheet2) get new border styles
#    (top+bottom thin, and top+bottom+right thin) reusing the workbook's
#    existing border definitions (borderId 4 and 5).
#  - Every "fedcore" label in the mode/approach header row becomes "approach".
#  - The stray empty inline-string cell G5 on the computational_comparison
#    sheet is removed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- build the two new border styles on sheet1, cells C1 (top+bottom) and
#     D1 (top+right+bottom) -------------------------------------------------
$c1 = $ws1.Cells.Item(1, 3)
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1
$c1.Borders.Item(9).LineStyle = 1

$d1 = $ws1.Cells.Item(1, 4)
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1
$d1.Borders.Item(10).LineStyle = 1
$d1.Borders.Item(9).LineStyle = 1

# --- reuse those two freshly built styles (via copy/paste of formats only)
#     for the matching cells on sheet2, so no redundant style-table entries
#     are produced -------------------------------------------------------
$c1.Copy()
$ws2.Cells.Item(1, 3).PasteSpecial(-4122)

$d1.Copy()
$ws2.Cells.Item(1, 4).PasteSpecial(-4122)

$c1.Copy()
$ws2.Cells.Item(1, 6).PasteSpecial(-4122)

$d1.Copy()
$ws2.Cells.Item(1, 7).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- rename every "fedcore" header label to "approach" ----------------------
$ws1.Cells.Item(2, 3).Value = "approach"

$ws2.Cells.Item(2, 3).Value = "approach"
$ws2.Cells.Item(2, 6).Value = "approach"

# --- drop the stray empty inline-string cell G5 on sheet2 -------------------
$ws2.Cells.Item(5, 7).ClearContents()
